# Draft for UseCase TipPrediction and RateDetection
#
# Resizes the small results table ("Tabelle 4") on slides 9-14 (the
# "Eingabe" header cell becomes "Trainingsdaten", and the table/columns
# are widened accordingly). Slide 10 additionally needs its row heights
# adjusted and the caption roundRect shape moved/shrunk.

$EMU_PER_POINT = 12700

$p = $ppt.ActivePresentation

function Set-TableHeaderAndGeometry {
    param(
        [int]$SlideIndex,
        [int]$OffX,
        [int]$OffY,
        [int]$ExtCx,
        [int]$ExtCy,
        [int[]]$ColWidths,
        [int[]]$RowHeights
    )

    $slide = $p.Slides.Item($SlideIndex)
    $shp = $slide.Shapes.Item("Tabelle 4")
    $tbl = $shp.Table

    # "Eingabe" -> "Trainingsdaten" in the first header cell, without
    # touching the following soft line-break / "( in Tausend)" run.
    $cell = $tbl.Cell(1, 1)
    $cellRange = $cell.Shape.TextFrame.TextRange
    $headerWord = $cellRange.Characters(1, 7)
    $headerWord.Text = "Trainingsdaten"

    # Column widths (EMU -> points).
    for ($c = 1; $c -le $ColWidths.Length; $c++) {
        $tbl.Columns.Item($c).Width = $ColWidths[$c - 1] / $EMU_PER_POINT
    }

    # Row heights, when this slide's table needs them touched.
    if ($RowHeights) {
        for ($r = 1; $r -le $RowHeights.Length; $r++) {
            $tbl.Rows.Item($r).Height = $RowHeights[$r - 1] / $EMU_PER_POINT
        }
    }

    # Final frame position/size (applied last so it wins over any
    # auto-layout nudging triggered by the column/row edits above).
    $shp.Left = $OffX / $EMU_PER_POINT
    $shp.Top = $OffY / $EMU_PER_POINT
    $shp.Width = $ExtCx / $EMU_PER_POINT
    $shp.Height = $ExtCy / $EMU_PER_POINT
}

# Slide 9 - Trinkgeldvorhersage (Regression)
Set-TableHeaderAndGeometry -SlideIndex 9 `
    -OffX 35496 -OffY 1902653 -ExtCx 3539093 -ExtCy 2171700 `
    -ColWidths @(1327159, 958505, 1253429)

# Slide 10 - Ratenerkennung (Multiklassen Klassifizierung)
Set-TableHeaderAndGeometry -SlideIndex 10 `
    -OffX 179512 -OffY 1902653 -ExtCx 3528392 -ExtCy 2312493 `
    -ColWidths @(1323146, 955608, 1249638) `
    -RowHeights @(552273, 243000, 243000, 243000, 243000, 243000, 243000, 243000)

$slide10 = $p.Slides.Item(10)
$caption = $slide10.Shapes.Item("Rechteck: abgerundete Ecken 6")
$caption.Top = 4299942 / $EMU_PER_POINT
$caption.Height = 288032 / $EMU_PER_POINT

# Slide 11 - Passagieraufkommen (Multiklassen-Klassifizierung)
Set-TableHeaderAndGeometry -SlideIndex 11 `
    -OffX 107504 -OffY 1902653 -ExtCx 3600399 -ExtCy 1668780 `
    -ColWidths @(1350148, 1064575, 1185676)

# Slide 12
Set-TableHeaderAndGeometry -SlideIndex 12 `
    -OffX 107504 -OffY 1902653 -ExtCx 3600400 -ExtCy 2171700 `
    -ColWidths @(1350149, 975110, 1275141)

# Slide 13
Set-TableHeaderAndGeometry -SlideIndex 13 `
    -OffX 262222 -OffY 1902653 -ExtCx 3589698 -ExtCy 2171700 `
    -ColWidths @(1346136, 972211, 1271351)

# Slide 14
Set-TableHeaderAndGeometry -SlideIndex 14 `
    -OffX 262222 -OffY 1902653 -ExtCx 3517690 -ExtCy 2171700 `
    -ColWidths @(1319133, 952709, 1245848)
